$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.648.71"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "'1.590.37"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'211.14"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").Value = "'0.509"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "'0.0615"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "'19.63"
$ws.Range("E10").Value = "  -3.89%  "
$ws.Range("D11").Value = "'0.0833"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "'1.811.06"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").Value = "'1.591.71"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "'0.523"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").Value = "'64.89"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'26.616.34"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "'0.0₃0728"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'208.48"
$ws.Range("E19").Value = "  -3.94%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.00"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  -3.10%  "
$ws.Range("E23").Value = "  -3.76%  "
$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "'146.81"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'7.26"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  -3.74%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "'0.0507"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").Value = "'0.664"
$ws.Range("E33").Value = "  +20.30%  "
$ws.Range("D34").Value = "'2.91"
$ws.Range("E34").Value = "  -2.77%  "
$ws.Range("D35").Value = "'1.308.63"
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").Value = "'1.49"
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("D38").Value = "'0.0172"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "'0.830"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "'0.793"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("D42").Value = "'5.37"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").Value = "'2.16"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").Value = "'62.87"
$ws.Range("E44").Value = "  -4.04%  "
$ws.Range("D45").Value = "'1.724.33"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "'89.64"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").Value = "'0.837"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").Value = "'0.0980"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "'7.54"
$ws.Range("E51").Value = "  -0.85%  "
